$d = $word.ActiveDocument

# Remove the transaction-data paragraph (2nd paragraph) entirely, including its
# paragraph mark, so "Apriori Analysis - Transaction 1" is immediately followed
# by the support/itemsets table paragraph.
$d.Paragraphs.Item(2).Range.Delete()

# Fix up the itemset tuple orderings and the antecedent/consequent table rows
# (and their dependent confidence/lift/zhangs_metric rows) to match the
# corrected notebook output.
$d.Content.Find.Execute('23  0.105263     (lettuce, cereal)', $true, $false, $false, $false, $false, $true, 1, $false, '23  0.105263     (cereal, lettuce)', 2) | Out-Null
$d.Content.Find.Execute('25  0.105263      (chips, lettuce)', $true, $false, $false, $false, $false, $true, 1, $false, '25  0.105263      (lettuce, chips)', 2) | Out-Null
$d.Content.Find.Execute('26  0.105263  (ice_cream, lettuce)', $true, $false, $false, $false, $false, $true, 1, $false, '26  0.105263  (lettuce, ice_cream)', 2) | Out-Null
$d.Content.Find.Execute('27  0.105263    (shampoo, lettuce)', $true, $false, $false, $false, $false, $true, 1, $false, '27  0.105263    (lettuce, shampoo)', 2) | Out-Null
$d.Content.Find.Execute('0    (lettuce)     (cereal)            0.263158            0.210526  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '0     (cereal)    (lettuce)            0.210526            0.263158  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('1     (cereal)    (lettuce)            0.210526            0.263158  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '1    (lettuce)     (cereal)            0.263158            0.210526  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('4      (chips)    (lettuce)            0.210526            0.263158  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '4    (lettuce)      (chips)            0.263158            0.210526  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('5    (lettuce)      (chips)            0.263158            0.210526  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '5      (chips)    (lettuce)            0.210526            0.263158  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('6  (ice_cream)    (lettuce)            0.157895            0.263158  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '6    (lettuce)  (ice_cream)            0.263158            0.157895  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('7    (lettuce)  (ice_cream)            0.263158            0.157895  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '7  (ice_cream)    (lettuce)            0.157895            0.263158  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('8    (shampoo)    (lettuce)            0.157895            0.263158  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '8    (lettuce)    (shampoo)            0.263158            0.157895  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('9    (lettuce)    (shampoo)            0.263158            0.157895  0.105263   ', $true, $false, $false, $false, $false, $true, 1, $false, '9    (shampoo)    (lettuce)            0.157895            0.263158  0.105263   ', 2) | Out-Null
$d.Content.Find.Execute('0    0.400000  1.900000               1.0  0.049861    1.315789   ', $true, $false, $false, $false, $false, $true, 1, $false, '0    0.500000  1.900000               1.0  0.049861    1.473684   ', 2) | Out-Null
$d.Content.Find.Execute('1    0.500000  1.900000               1.0  0.049861    1.473684   ', $true, $false, $false, $false, $false, $true, 1, $false, '1    0.400000  1.900000               1.0  0.049861    1.315789   ', 2) | Out-Null
$d.Content.Find.Execute('4    0.500000  1.900000               1.0  0.049861    1.473684   ', $true, $false, $false, $false, $false, $true, 1, $false, '4    0.400000  1.900000               1.0  0.049861    1.315789   ', 2) | Out-Null
$d.Content.Find.Execute('5    0.400000  1.900000               1.0  0.049861    1.315789   ', $true, $false, $false, $false, $false, $true, 1, $false, '5    0.500000  1.900000               1.0  0.049861    1.473684   ', 2) | Out-Null
$d.Content.Find.Execute('6    0.666667  2.533333               1.0  0.063712    2.210526   ', $true, $false, $false, $false, $false, $true, 1, $false, '6    0.400000  2.533333               1.0  0.063712    1.403509   ', 2) | Out-Null
$d.Content.Find.Execute('7    0.400000  2.533333               1.0  0.063712    1.403509   ', $true, $false, $false, $false, $false, $true, 1, $false, '7    0.666667  2.533333               1.0  0.063712    2.210526   ', 2) | Out-Null
$d.Content.Find.Execute('8    0.666667  2.533333               1.0  0.063712    2.210526   ', $true, $false, $false, $false, $false, $true, 1, $false, '8    0.400000  2.533333               1.0  0.063712    1.403509   ', 2) | Out-Null
$d.Content.Find.Execute('9    0.400000  2.533333               1.0  0.063712    1.403509   ', $true, $false, $false, $false, $false, $true, 1, $false, '9    0.666667  2.533333               1.0  0.063712    2.210526   ', 2) | Out-Null
$d.Content.Find.Execute('0       0.642857  0.285714   0.240000    0.450000  ', $true, $false, $false, $false, $false, $true, 1, $false, '0       0.600000  0.285714   0.321429    0.450000  ', 2) | Out-Null
$d.Content.Find.Execute('1       0.600000  0.285714   0.321429    0.450000  ', $true, $false, $false, $false, $false, $true, 1, $false, '1       0.642857  0.285714   0.240000    0.450000  ', 2) | Out-Null
$d.Content.Find.Execute('4       0.600000  0.285714   0.321429    0.450000  ', $true, $false, $false, $false, $false, $true, 1, $false, '4       0.642857  0.285714   0.240000    0.450000  ', 2) | Out-Null
$d.Content.Find.Execute('5       0.642857  0.285714   0.240000    0.450000  ', $true, $false, $false, $false, $false, $true, 1, $false, '5       0.600000  0.285714   0.321429    0.450000  ', 2) | Out-Null
$d.Content.Find.Execute('6       0.718750  0.333333   0.547619    0.533333  ', $true, $false, $false, $false, $false, $true, 1, $false, '6       0.821429  0.333333   0.287500    0.533333  ', 2) | Out-Null
$d.Content.Find.Execute('7       0.821429  0.333333   0.287500    0.533333  ', $true, $false, $false, $false, $false, $true, 1, $false, '7       0.718750  0.333333   0.547619    0.533333  ', 2) | Out-Null
$d.Content.Find.Execute('8       0.718750  0.333333   0.547619    0.533333  ', $true, $false, $false, $false, $false, $true, 1, $false, '8       0.821429  0.333333   0.287500    0.533333  ', 2) | Out-Null
$d.Content.Find.Execute('9       0.821429  0.333333   0.287500    0.533333  ', $true, $false, $false, $false, $false, $true, 1, $false, '9       0.718750  0.333333   0.547619    0.533333  ', 2) | Out-Null
